# Fix inconsistent latitude/longitude and kiosk counts for "Centre A Ben Mansour"
# (rows 2-4, id 10001), and correct the number_of_kiosks value (column L) for
# every registration-center row to reflect how many kiosk rows belong to
# each center.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Latitude (column G): rows 3 and 4 had slightly-off values (35.52.., 36.52..)
# instead of matching row 2's 34.52.. -- align them.
$ws.Range("G3").Value = $ws.Range("G2").Value()
$ws.Range("G4").Value = $ws.Range("G2").Value()

# --- Longitude (column H): rows 3 and 4 referenced two near-duplicate text
# values (" -6.453276", " -6.453277" with a leading non-breaking space)
# instead of row 2's " -6.453275". Copy H2 straight onto H3/H4 so they end
# up pointing at the exact same shared-string text/type as H2 (rather than
# re-typing the text, which Excel would reparse as a number) -- the two
# stray strings then drop out of the shared-string table entirely.
$ws.Range("H2").Copy($ws.Range("H3"))
$ws.Range("H2").Copy($ws.Range("H4"))

# --- number_of_kiosks (column L): was hard-coded to 1 everywhere; correct it
# to the real kiosk count per center (3 kiosks for center 10001, 2 kiosks for
# every other center).
$ws.Range("L2:L4").Value = 3
$ws.Range("L5:L46").Value = 2
